## Update with Correct Forecast output
## Restructure the single-sheet sales/PO comparison workbook into four
## sheets: "Sales vs PO" (original data with a new Order Week column),
## "Weekly Growth" (PO qty + week-over-week growth%), "Volume Insights"
## (aggregate PO stats) and "Prediction Info" (next week forecast).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original sheet and insert the new "Order Week" column.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Insert a new column C (pushes the old PO_Requested_Qty column to D).
$ws1.Columns.Item(3).Insert()
$ws1.Range("C1").Value = "Order Week"

$lastRow = 18
for ($r = 2; $r -le $lastRow; $r++) {
    $origDs = $ws1.Cells.Item($r, 1).Value2

    # New "Order Week" column keeps the original ds date value/format.
    $ws1.Cells.Item($r, 3).Value = $origDs
    $ws1.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # The ds column itself shifts forward one week (+6 days from the
    # previous value, i.e. the following week's date).
    $ws1.Cells.Item($r, 1).Value = $origDs + 6

    # PO_Requested_Qty now lives on the "Weekly Growth" sheet; zero this
    # column out here.
    $ws1.Cells.Item($r, 4).Value = 0
}

# ---------------------------------------------------------------------
# 2. "Weekly Growth" sheet.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$wgDates  = @(45544, 45572, 45579, 45586, 45593)
$wgQty    = @(500, 630, 200, 10, 750)
$wgGrowth = @(0, 26, -68.25396825396825, -95, 7400)

for ($i = 0; $i -lt $wgDates.Length; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 1).Value = $wgDates[$i]
    $ws2.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws2.Cells.Item($r, 2).Value = $wgQty[$i]
    $ws2.Cells.Item($r, 3).Value = $wgGrowth[$i]
}

# ---------------------------------------------------------------------
# 3. "Volume Insights" sheet.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value = 2090
$ws3.Range("B2").Value = 418
$ws3.Range("C2").Value = 750
$ws3.Range("D2").Value = 10

# ---------------------------------------------------------------------
# 4. "Prediction Info" sheet.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Range("A2").Value = 381.9999999999999

Write-Output "done"
